# "Made delete work and started on sorts"
#
# The sheet used to hold a tiny 2-row placeholder table in A1:F3. It is
# replaced here by a 5-row table in A1:H6 that also gains two new leading
# data columns ("Unnamed: 0.1" / "Unnamed: 0") - the shape pandas leaves
# behind after an index reset/delete - with the rows re-sorted.
#
# NOTE: cells B1:F1 and A2:A3 already carry the bold/centered/bordered
# header-ish style (style index 1) in the original workbook, and simply
# overwriting .Value on a cell does not disturb its existing style/format.
# So we only need to explicitly copy that format onto the brand-new cells
# that fall outside the old A1:F3 range (G1, H1, A4:A6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ------------------------------------------------
# B1/C1 are brand-new text headers (they reuse the style already sitting
# on those cells); D1:H1 keep the old numeric 0..4 header, shifted two
# columns to the right.
$ws.Cells.Item(1, 2).Value = "Unnamed: 0.1"
$ws.Cells.Item(1, 3).Value = "Unnamed: 0"
$ws.Cells.Item(1, 4).Value = 0
$ws.Cells.Item(1, 5).Value = 1
$ws.Cells.Item(1, 6).Value = 2
$ws.Cells.Item(1, 7).Value = 3
$ws.Cells.Item(1, 8).Value = 4

# G1/H1 are new cells (outside the old A1:F3 range) - give them the same
# format as the rest of the header row by copying it from an already
# styled header cell.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null

# --- Data rows (rows 2-6) ----------------------------------------------
# Columns: A=index, B="Unnamed: 0.1", C="Unnamed: 0", D=col0, E=Name,
#          F=amount, G=amount/10, H=id
$rows = @(
    @{ Row = 2; A = 0; B = 0; C = 1; D = 5; E = "Transaction 1"; F = 100; G = 10; H = 1 },
    @{ Row = 3; A = 1; B = 4; C = 5; D = 1; E = "Transaction 5"; F = 500; G = 50; H = 5 },
    @{ Row = 4; A = 2; B = 3; C = 4; D = 2; E = "Transaction 4"; F = 400; G = 40; H = 4 },
    @{ Row = 5; A = 3; B = 2; C = 3; D = 3; E = "Transaction 3"; F = 300; G = 30; H = 3 },
    @{ Row = 6; A = 4; B = 1; C = 2; D = 4; E = "Transaction 2"; F = 200; G = 20; H = 2 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}

# A4:A6 are new index cells (outside the old A1:F3 range) - copy the same
# bold/centered/bordered format already on A2/A3 down onto them.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A6").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Select() | Out-Null
